$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.928.65'
$ws.Range("E2").Value = '  -0.87%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.355.54'
$ws.Range("E3").Value = '  -0.52%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.04'
$ws.Range("E5").Value = '  -1.15%  '

$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.669'
$ws.Range("E6").Value = '  -4.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.42'
$ws.Range("E7").Value = '  -1.93%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.604'
$ws.Range("E9").Value = '  +0.59%  '

$ws.Range("E10").Value = '  -2.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.23'
$ws.Range("E11").Value = '  +2.46%  '

$ws.Range("E12").Value = '  +5.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.34'
$ws.Range("E13").Value = '  -2.27%  '

$ws.Range("E14").Value = '  -0.97%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.705.49'

$ws.Range("E16").Value = '  -3.76%  '

$ws.Range("E17").Value = '  -0.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.352.31'
$ws.Range("E18").Value = '  -0.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.847.50'
$ws.Range("E19").Value = '  -0.97%  '

$ws.Range("E20").Value = '  -0.75%  '

$ws.Range("E21").Value = '  +0.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.74'
$ws.Range("E22").Value = '  -1.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '257.17'
$ws.Range("E23").Value = '  -0.32%  '

$ws.Range("E24").Value = '  +15.15%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("E27").Value = '  -2.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.65'
$ws.Range("E28").Value = '  -1.80%  '

$ws.Range("E29").Value = '  +0.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.73'
$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '177.43'
$ws.Range("E31").Value = '  +1.41%  '

$ws.Range("E32").Value = '  -0.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.136'
$ws.Range("E33").Value = '  +0.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0758'
$ws.Range("E34").Value = '  -0.09%  '

$ws.Range("E35").Value = '  -3.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.52'
$ws.Range("E36").Value = '  +2.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.81'
$ws.Range("E37").Value = '  -2.86%  '

$ws.Range("E38").Value = '  -1.94%  '

$ws.Range("E39").Value = '  -4.05%  '

$ws.Range("E40").Value = '  +0.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '68.36'
$ws.Range("E41").Value = '  +27.65%  '

$ws.Range("E42").Value = '  +10.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.14'
$ws.Range("E43").Value = '  +15.60%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.38'
$ws.Range("E44").Value = '  +2.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.203'
$ws.Range("E45").Value = '  +2.65%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.09'
$ws.Range("E46").Value = '  -1.19%  '

$ws.Range("E47").Value = '  +0.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.50'
$ws.Range("E48").Value = '  +0.49%  '

$ws.Range("E49").Value = '  -0.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '99.73'
$ws.Range("E50").Value = '  -1.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.16'
$ws.Range("E51").Value = '  -3.68%  '
